$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: append plain text to the very end of a paragraph (just before its
# trailing paragraph mark), picking up whatever formatting is already there.
# ---------------------------------------------------------------------------
function Append-PlainRun($para, [string]$text) {
    $pos = $para.Range.End - 1
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
}

# ---------------------------------------------------------------------------
# Helper: append text to the end of a paragraph and then explicitly stamp the
# newly inserted span with Roboto / color 424242 formatting.
# ---------------------------------------------------------------------------
function Append-FormattedRun($para, [string]$text) {
    $startPos = $para.Range.End - 1
    $r = $d.Range($startPos, $startPos)
    $r.InsertAfter($text)
    $endPos = $startPos + $text.Length
    $rf = $d.Range($startPos, $endPos)
    $rf.Font.Name = "Roboto"
    $rf.Font.Color = 4342338
}

# ---------------------------------------------------------------------------
# Locate the "Re-sized all images to be 800px" bullet paragraph, and the
# (empty) bullet paragraph that immediately follows it.
# ---------------------------------------------------------------------------
$resizedPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd()
    if ($txt -eq "Re-sized all images to be 800px") {
        $resizedPara = $p
    }
}

if ($resizedPara -eq $null) {
    Write-Host "Could not locate 'Re-sized all images to be 800px' paragraph"
} else {
    # Extend "Re-sized all images to be 800px" with the trailing sentence.
    Append-PlainRun $resizedPara " and"
    Append-PlainRun $resizedPara " load faster and consume less cellular data"
    Append-PlainRun $resizedPara "."

    # The next bullet (currently empty) gets the "Only used ..." sentence in
    # Roboto / #424242.
    $nextPara = $resizedPara.Next()
    Append-FormattedRun $nextPara "Only used "
    Append-FormattedRun $nextPara "critical JS/CSS inline and deferring all non-critical JS/styles"
    Append-FormattedRun $nextPara " because r"
    Append-FormattedRun $nextPara "esources are blocking the first paint of "
    Append-FormattedRun $nextPara "my"
    Append-FormattedRun $nextPara " page."
}
